# Regenerate save_data column G ("K") values for wittgren_nick.xlsx
# The underlying source switched from using a "Strike#" metric to true
# strikeout counts ("K"), so the raw per-appearance values in column G
# (rows 2-34) are rewritten with the newly computed strikeout totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 3
    6  = 0
    7  = 1
    8  = 1
    9  = 0
    10 = 2
    11 = 0
    12 = 0
    13 = 0
    14 = 0
    15 = 0
    16 = 0
    17 = 1
    18 = 0
    19 = 0
    20 = 0
    21 = 1
    22 = 1
    23 = 1
    24 = 1
    25 = 0
    26 = 1
    27 = 0
    28 = 0
    29 = 1
    30 = 0
    31 = 0
    32 = 1
    33 = 2
    34 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
